$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "37.806.60"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "2.087.63"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "234.73"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").Value = "58.85"
$ws.Range("E7").Value = "  +3.26%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.393"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").Value = "0.0792"
$ws.Range("E10").Value = "  +2.10%  "
$ws.Range("D11").Value = "0.106"
$ws.Range("E11").Value = "  +2.81%  "
$ws.Range("D12").Value = "2.394.16"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("D13").Value = "14.78"
$ws.Range("E13").Value = "  +3.34%  "
$ws.Range("D14").Value = "21.27"
$ws.Range("E14").Value = "  +3.21%  "
$ws.Range("D15").Value = "0.770"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").Value = "5.32"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").Value = "2.083.93"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").Value = "37.720.99"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").Value = "6.22"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "71.33"
$ws.Range("E20").Value = "  +2.65%  "
$ws.Range("E21").Value = "  +2.41%  "
$ws.Range("D22").Value = "228.78"
$ws.Range("E22").Value = "  +1.17%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("D25").Value = "2.40"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("D26").Value = "170.22"
$ws.Range("E27").Value = "  +5.09%  "
$ws.Range("D28").Value = "9.01"
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "19.53"
$ws.Range("E30").Value = "  +2.42%  "
$ws.Range("D31").Value = "0.121"
$ws.Range("E31").Value = "  +2.51%  "
$ws.Range("E32").Value = "  +3.19%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "4.72"
$ws.Range("E33").Value = "  +4.31%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.0633"
$ws.Range("E34").Value = "  +2.69%  "
$ws.Range("E35").Value = "  +1.79%  "
$ws.Range("D36").Value = "3.47"
$ws.Range("E36").Value = "  +3.31%  "
$ws.Range("E37").Value = "  +3.25%  "
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("E39").Value = "  -3.94%  "
$ws.Range("D40").Value = "0.0995"
$ws.Range("E40").Value = "  +4.60%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "98.88"
$ws.Range("E41").Value = "  +2.16%  "
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "2.94"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0216"
$ws.Range("E43").Value = "  +1.55%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "4.37"
$ws.Range("E44").Value = "  +4.56%  "
$ws.Range("D45").Value = "1.466.85"
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("D47").Value = "1.08"
$ws.Range("E47").Value = "  +4.42%  "
$ws.Range("D48").Value = "16.07"
$ws.Range("E48").Value = "  +5.77%  "
$ws.Range("D49").Value = "7.41"
$ws.Range("E49").Value = "  +2.77%  "
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("D51").Value = "2.278.73"
$ws.Range("E51").Value = "  +0.85%  "
